# Citation-check update: swap the placeholder reference tags for the
# real in-text citations, scoped per-paragraph so that the same
# placeholder text (e.g. "Ref-A1B2C3") can resolve to different
# citations depending on which paragraph it appears in.

$d = $word.ActiveDocument

# Paragraph 2: "The use of social media and educational campaigns..."
# All four distinct placeholders collapse onto the same citation.
$p2 = $d.Paragraphs(2).Range
$p2.Find.Execute("Ref-A1B2C3", $true, $false, $false, $false, $false, $true, 1, $false, "Ref-s232886", 2)
$p2 = $d.Paragraphs(2).Range
$p2.Find.Execute("Ref-D4E5F6", $true, $false, $false, $false, $false, $true, 1, $false, "Ref-s232886", 2)
$p2 = $d.Paragraphs(2).Range
$p2.Find.Execute("Ref-G7H8I9", $true, $false, $false, $false, $false, $true, 1, $false, "Ref-s232886", 2)
$p2 = $d.Paragraphs(2).Range
$p2.Find.Execute("Ref-J0K1L2", $true, $false, $false, $false, $false, $true, 1, $false, "Ref-s232886", 2)

# Paragraph 3: "According to the study of Mim and Jameelah..."
# All three placeholders become the author-date citation.
$p3 = $d.Paragraphs(3).Range
$p3.Find.Execute("Ref-A1B2C3", $true, $false, $false, $false, $false, $true, 1, $false, "Brown & Garcia, 2018", 2)
$p3 = $d.Paragraphs(3).Range
$p3.Find.Execute("Ref-D4E5F6", $true, $false, $false, $false, $false, $true, 1, $false, "Brown & Garcia, 2018", 2)
$p3 = $d.Paragraphs(3).Range
$p3.Find.Execute("Ref-G7H8I9", $true, $false, $false, $false, $false, $true, 1, $false, "Brown & Garcia, 2018", 2)

# Paragraph 4: "There are genetic factors that could affect the infant..."
$p4 = $d.Paragraphs(4).Range
$p4.Find.Execute("Ref-G7H8I9", $true, $false, $false, $false, $false, $true, 1, $false, "Ref-s084922", 2)
$p4 = $d.Paragraphs(4).Range
$p4.Find.Execute("Ref-J1K2L3", $true, $false, $false, $false, $false, $true, 1, $false, "Ref-s084922", 2)

# Paragraph 5: "Moreover, 40% of cases of SIDS..."
$p5 = $d.Paragraphs(5).Range
$p5.Find.Execute("Ref-AB1CD2", $true, $false, $false, $false, $false, $true, 1, $false, "Johnson 45", 2)
$p5 = $d.Paragraphs(5).Range
$p5.Find.Execute("Ref-EF3GH4", $true, $false, $false, $false, $false, $true, 1, $false, "Johnson 45", 2)
$p5 = $d.Paragraphs(5).Range
$p5.Find.Execute("Ref-IJ5KL6", $true, $false, $false, $false, $false, $true, 1, $false, "Johnson 45", 2)
